# Fixed update to excel issue
#
# 1) Rename the "Requested quantity" header on the two existing sheets.
# 2) Add a new "PO Forecast" sheet (as the 3rd / last tab) with forecast data.

$wb = $excel.ActiveWorkbook

# --- 1) Rename headers -----------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the "PO Forecast" sheet ----------------------------------------
$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows
$data = @(
  @(45375.99999999999, 18, 2.436238119479508, 31.62961692170215),
  @(45382.99999999999, 17, 3.445757120390903, 31.41918233988922),
  @(45396.99999999999, 17, 3.471544585147041, 31.31779855440318),
  @(45424.99999999999, 17, 1.611155335884578, 30.52357970209304),
  @(45431.99999999999, 16, 2.792614597837289, 31.86361268155539),
  @(45494.99999999999, 15, 1.257641666324014, 30.84108228763887),
  @(45557.99999999999, 14, -1.356476786434347, 27.72249343039077),
  @(45578.99999999999, 13, -1.058855895478787, 27.079944765331),
  @(45599.99999999999, 13, -0.827558350650939, 27.12325272481121),
  @(45606.99999999999, 13, -1.53566134098037, 27.80842136336634),
  @(45613.99999999999, 13, -1.394010755514835, 26.05137353072757),
  @(45620.99999999999, 12, -2.117196224577669, 25.71963050089251),
  @(45627.99999999999, 12, -2.378199405314541, 25.85273435234149),
  @(45634.99999999999, 12, -2.389413685317954, 26.70634092903932),
  @(45641.99999999999, 12, -1.662846231613039, 26.02952623390584),
  @(45648.99999999999, 12, -2.775679403568247, 26.05262056171091),
  @(45655.99999999999, 12, -3.357563164644826, 25.31419394213078)
)

$row = 2
foreach ($r in $data) {
  $wsForecast.Cells.Item($row, 1).Value = $r[0]
  $wsForecast.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
  $wsForecast.Cells.Item($row, 2).Value = $r[1]
  $wsForecast.Cells.Item($row, 3).Value = $r[2]
  $wsForecast.Cells.Item($row, 4).Value = $r[3]
  $row = $row + 1
}

# Move the new sheet to the end of the tab strip (after "Monthly Trend")
$wsForecast.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
